# "revise examples in Appendix A"
#
# Each worksheet in this workbook holds two stacked copies of the same
# example block, separated by a few blank spacer rows. The revision tightens
# the spacing between the two blocks on every sheet by removing extra blank
# spacer row(s), which shifts the second block up. It also moves the active
# tab/selection from "CONNECT Examples" to "PUBLISH example".

$wb = $excel.ActiveWorkbook

$wsConnect = $wb.Worksheets.Item("CONNECT Examples")
$wsPublish = $wb.Worksheets.Item("PUBLISH example")
$wsSubscribe = $wb.Worksheets.Item("zulu SUBSCRIBE")

# --- "CONNECT Examples": 2 blank rows -> 1 blank row before the ruled row ---
$wsConnect.Rows.Item(11).Delete() | Out-Null
$wsConnect.Range("G9").Select() | Out-Null

# --- "PUBLISH example": 3 blank rows -> 1 blank row before the ruled row ---
$wsPublish.Rows.Item(13).Delete() | Out-Null
$wsPublish.Rows.Item(13).Delete() | Out-Null
$wsPublish.Range("G15").Select() | Out-Null

# --- "zulu SUBSCRIBE": 2 blank rows -> 1 blank row before the ruled row ---
$wsSubscribe.Rows.Item(14).Delete() | Out-Null
$wsSubscribe.Range("K15").Select() | Out-Null

# Make "PUBLISH example" the active/visible tab.
$wsPublish.Activate() | Out-Null
